$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3399353333333333
$ws.Range("H2").Value = 1.019806
$ws.Range("I2").Value = 0.09929991924017606
$ws.Range("J2").Value = 0.09929991924017606
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.023286
$ws.Range("N2").Value = 0.069858
$ws.Range("O2").Value = 0.009310710475795457
$ws.Range("P2").Value = 0.009310710475795458
$ws.Range("Q2").Value = 0.007915734172
$ws.Range("R2").Value = 0.071241607548
$ws.Range("S2").Value = 0.00092455279831515
$ws.Range("T2").Value = 0.0009245527983151502
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3399353333333333
$ws.Range("H3").Value = 1.019806
$ws.Range("I3").Value = 0.09929991924017606
$ws.Range("J3").Value = 0.09929991924017606
$ws.Range("O3").Value = 0.05314667307834813
$ws.Range("P3").Value = 0.05314667307834814
$ws.Range("Q3").Value = 0.0451839778831111
$ws.Range("R3").Value = 0.406655800948
$ws.Range("S3").Value = 0.005277460344564009
$ws.Range("T3").Value = 0.00527746034456401
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3399353333333333
$ws.Range("H4").Value = 1.019806
$ws.Range("I4").Value = 0.09929991924017606
$ws.Range("J4").Value = 0.09929991924017606
$ws.Range("M4").Value = 2.344785333333334
$ws.Range("N4").Value = 7.034356000000001
$ws.Range("O4").Value = 0.9375426164458565
$ws.Range("P4").Value = 0.9375426164458565
$ws.Range("Q4").Value = 0.7970753838817779
$ws.Range("R4").Value = 7.173678454936001
$ws.Range("S4").Value = 0.0930979060972969
$ws.Range("T4").Value = 0.0930979060972969
$ws.Range("I5").Value = 0.4094685684206303
$ws.Range("J5").Value = 0.4094685684206303
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.023286
$ws.Range("N5").Value = 0.069858
$ws.Range("O5").Value = 0.009310710475795457
$ws.Range("P5").Value = 0.009310710475795458
$ws.Range("Q5").Value = 0.03264095644999999
$ws.Range("R5").Value = 0.29376860805
$ws.Range("S5").Value = 0.003812443289502931
$ws.Range("T5").Value = 0.003812443289502932
$ws.Range("I6").Value = 0.4094685684206303
$ws.Range("J6").Value = 0.4094685684206303
$ws.Range("O6").Value = 0.05314667307834813
$ws.Range("P6").Value = 0.05314667307834814
$ws.Range("S6").Value = 0.02176189214171046
$ws.Range("T6").Value = 0.02176189214171047
$ws.Range("I7").Value = 0.4094685684206303
$ws.Range("J7").Value = 0.4094685684206303
$ws.Range("M7").Value = 2.344785333333334
$ws.Range("N7").Value = 7.034356000000001
$ws.Range("O7").Value = 0.9375426164458565
$ws.Range("P7").Value = 0.9375426164458565
$ws.Range("Q7").Value = 3.286783301122222
$ws.Range("R7").Value = 29.5810497101
$ws.Range("S7").Value = 0.3838942329894169
$ws.Range("T7").Value = 0.3838942329894169
$ws.Range("G8").Value = 1.681642333333333
$ws.Range("H8").Value = 5.044927
$ws.Range("I8").Value = 0.4912315123391937
$ws.Range("J8").Value = 0.4912315123391937
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.023286
$ws.Range("N8").Value = 0.069858
$ws.Range("O8").Value = 0.009310710475795457
$ws.Range("P8").Value = 0.009310710475795458
$ws.Range("Q8").Value = 0.039158723374
$ws.Range("R8").Value = 0.352428510366
$ws.Range("S8").Value = 0.004573714387977375
$ws.Range("T8").Value = 0.004573714387977376
$ws.Range("G9").Value = 1.681642333333333
$ws.Range("H9").Value = 5.044927
$ws.Range("I9").Value = 0.4912315123391937
$ws.Range("J9").Value = 0.4912315123391937
$ws.Range("O9").Value = 0.05314667307834813
$ws.Range("P9").Value = 0.05314667307834814
$ws.Range("Q9").Value = 0.2235227778517777
$ws.Range("R9").Value = 2.011705000666
$ws.Range("S9").Value = 0.02610732059207366
$ws.Range("T9").Value = 0.02610732059207366
$ws.Range("G10").Value = 1.681642333333333
$ws.Range("H10").Value = 5.044927
$ws.Range("I10").Value = 0.4912315123391937
$ws.Range("J10").Value = 0.4912315123391937
$ws.Range("M10").Value = 2.344785333333334
$ws.Range("N10").Value = 7.034356000000001
$ws.Range("O10").Value = 0.9375426164458565
$ws.Range("P10").Value = 0.9375426164458565
$ws.Range("Q10").Value = 3.943090279112445
$ws.Range("R10").Value = 35.48781251201201
$ws.Range("S10").Value = 0.4605504773591426
$ws.Range("T10").Value = 0.4605504773591426
